$wb = $excel.ActiveWorkbook

# --- Sheet: System Configuration0@0x0 ---
# Radio Name (Bluetooth advertisement name) changed for STU test
$wsConfig = $wb.Worksheets.Item("System Configuration0@0x0")
$wsConfig.Cells.Item(3, 5).Value = "BBCB4866"

# --- Sheet: Statistics@0x5 ---
$wsStats = $wb.Worksheets.Item("Statistics@0x5")
# Production Date
$wsStats.Cells.Item(7, 5).Value = "20191212"
# Batch Number STH
$wsStats.Cells.Item(8, 5).Value = "98"
# New row 9: extra value (e.g. device identifier / Bluetooth address for new naming)
$wsStats.Cells.Item(9, 5).Value = "9259266508322"

# --- Sheet: Calibration0@0x8 ---
$wsCal = $wb.Worksheets.Item("Calibration0@0x8")
# Acceleration X - K
$wsCal.Cells.Item(2, 5).Value = "0.0030518043786287308"
# Acceleration X - D
$wsCal.Cells.Item(3, 5).Value = "-99.06767272949219"
# Voltage Battery - D
$wsCal.Cells.Item(9, 5).Value = "-0.051661375910043716"
